$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "I7C8C9"
$ws.Range("C3").Value = "QRZEDJ"
$ws.Range("C4").Value = "KXAPBW"
$ws.Range("C5").Value = "WCBJT6"
$ws.Range("C6").Value = "HI7NMR"
$ws.Range("C7").Value = "UB7O8J"
$ws.Range("C8").Value = "ZX68WK"
$ws.Range("C9").Value = "3PQZYA"
$ws.Range("C10").Value = "30TCTL"
$ws.Range("C11").Value = "MR7C0R"
$ws.Range("C12").Value = "EJKKAU"
$ws.Range("C13").Value = "RAU62A"
$ws.Range("C14").Value = "21U4L4"
$ws.Range("C15").Value = "DF7LO8"
$ws.Range("C16").Value = "4A2GPD"
$ws.Range("C17").Value = "I6V0MG"
$ws.Range("C18").Value = "0WTJ1U"
$ws.Range("C19").Value = "1P1IEV"
$ws.Range("C20").Value = "7FIZDH"
$ws.Range("C21").Value = "KIDAXE"
$ws.Range("C22").Value = "O8QQA7"
$ws.Range("C23").Value = "DTENCZ"
$ws.Range("C24").Value = "167RZ1"
$ws.Range("C25").Value = "3D8DMU"
$ws.Range("C26").Value = "4MOV97"
$ws.Range("C27").Value = "PYG18V"
$ws.Range("C28").Value = "79XORS"
$ws.Range("C29").Value = "GBHZQQ"
$ws.Range("C30").Value = "EPY5UG"
$ws.Range("C31").Value = "GNGUL0"
$ws.Range("C32").Value = "FI21GO"
$ws.Range("C33").Value = "DC6MCO"
$ws.Range("C34").Value = "8V9WJD"
$ws.Range("C35").Value = "EKH9GV"
$ws.Range("C36").Value = "4YIHHU"
$ws.Range("C37").Value = "DIM0P2"
$ws.Range("C38").Value = "VN4N2U"
$ws.Range("C39").Value = "F11HA0"
$ws.Range("C40").Value = "IXVC0B"
$ws.Range("C41").Value = "LFKOGQ"
$ws.Range("C42").Value = "IDB46C"
$ws.Range("C43").Value = "AL76HK"
$ws.Range("C44").Value = "83Q0EZ"
$ws.Range("C45").Value = "FFN9LM"
$ws.Range("C46").Value = "4IHW1N"
$ws.Range("C47").Value = "985MIQ"
$ws.Range("C48").Value = "Z0R25N"
$ws.Range("C49").Value = "JXMUOY"
$ws.Range("C50").Value = "LA5JL7"
$ws.Range("C51").Value = "EYNKKC"
$ws.Range("C52").Value = "MFC3IP"
$ws.Range("C53").Value = "429NGL"
$ws.Range("C54").Value = "ZCHBEM"
$ws.Range("C55").Value = "6VHX8V"
$ws.Range("C56").Value = "W9XWGW"
$ws.Range("C57").Value = "WHR2H5"
$ws.Range("C58").Value = "UW4X4W"
$ws.Range("C59").Value = "VAVPG5"
$ws.Range("C60").Value = "43J8SU"
$ws.Range("C61").Value = "5AP4J3"
$ws.Range("C62").Value = "HI4VED"
$ws.Range("C63").Value = "3W745W"
$ws.Range("C64").Value = "VWT4T0"
$ws.Range("C65").Value = "2YUP88"
$ws.Range("C66").Value = "B9NRNA"
$ws.Range("C67").Value = "VMT8B9"
$ws.Range("C68").Value = "1WVEKD"
$ws.Range("C69").Value = "6M87G5"
$ws.Range("C70").Value = "KSBUSK"
$ws.Range("C71").Value = "CVH51T"
$ws.Range("C72").Value = "S6MYLX"
$ws.Range("C73").Value = "F0O75J"
$ws.Range("C74").Value = "X9LX0C"
$ws.Range("C75").Value = "7Q0YSZ"
$ws.Range("C76").Value = "SYQUK2"
$ws.Range("C77").Value = "ZVDMRP"
$ws.Range("C78").Value = "9R09ZG"
$ws.Range("C79").Value = "V90BL4"
$ws.Range("C80").Value = "Q8OH7J"
$ws.Range("C81").Value = "F047TI"
$ws.Range("C82").Value = "2RGAUS"
$ws.Range("C83").Value = "E2OVFQ"
$ws.Range("C84").Value = "5PEHPZ"
$ws.Range("C85").Value = "18QCJM"
$ws.Range("C86").Value = "3H5TLB"
$ws.Range("C87").Value = "1SB9LJ"
$ws.Range("C88").Value = "L18JAH"
$ws.Range("C89").Value = "W3HUP1"
$ws.Range("C90").Value = "6VHHY8"
$ws.Range("C91").Value = "533TTO"
$ws.Range("C92").Value = "BCM8UT"
$ws.Range("C93").Value = "JKBOE3"
$ws.Range("C94").Value = "RYNSWU"
$ws.Range("C95").Value = "JTFCG3"
$ws.Range("C96").Value = "2WJTMJ"
$ws.Range("C97").Value = "9VRG7C"
$ws.Range("C98").Value = "NQLLEG"
$ws.Range("C99").Value = "WPNGMJ"
$ws.Range("C100").Value = "5CSD1C"
$ws.Range("C101").Value = "3OA0R9"
$ws.Range("C102").Value = "ME2HBS"
$ws.Range("C103").Value = "UZBCI7"
$ws.Range("C104").Value = "6LYQA7"
$ws.Range("C105").Value = "TXYXQQ"
$ws.Range("C106").Value = "TIJXHV"
$ws.Range("C107").Value = "FETEMN"
$ws.Range("C108").Value = "2TUY6E"
$ws.Range("C109").Value = "DBN5GZ"
$ws.Range("C110").Value = "ORNQAJ"
$ws.Range("C111").Value = "N6T3ZB"
$ws.Range("C112").Value = "4JU5BY"
$ws.Range("C113").Value = "BNCZW8"
$ws.Range("C114").Value = "BIRZN7"
$ws.Range("C115").Value = "FZKA6J"
$ws.Range("C116").Value = "S3MLLR"
$ws.Range("C117").Value = "N6L3N2"
$ws.Range("C118").Value = "LTU9CP"
$ws.Range("C119").Value = "9B10UR"
$ws.Range("C120").Value = "R4YFBS"
$ws.Range("C121").Value = "K4XGQO"
$ws.Range("C122").Value = "H4KPOQ"
$ws.Range("C123").Value = "3Q5XK0"
$ws.Range("C124").Value = "OIU9ZK"
$ws.Range("C125").Value = "WGP0L1"
$ws.Range("C126").Value = "RI2PQ2"
$ws.Range("C127").Value = "VEIR9F"
$ws.Range("C128").Value = "Y3F6X6"
$ws.Range("C129").Value = "CLLOI3"
$ws.Range("C130").Value = "WRP1JM"
$ws.Range("C131").Value = "8EKNSD"
$ws.Range("C132").Value = "G7JHG6"
$ws.Range("C133").Value = "4Z4R7U"
$ws.Range("C134").Value = "G99B8Z"
$ws.Range("C135").Value = "LY0RKU"
$ws.Range("C136").Value = "T4VT8I"
$ws.Range("C137").Value = "GWKC62"
$ws.Range("C138").Value = "XPA9GI"
$ws.Range("C139").Value = "X8IO8N"
$ws.Range("C140").Value = "7I2MYY"
$ws.Range("C141").Value = "01ZC3F"
$ws.Range("C142").Value = "8QVE5U"
$ws.Range("C143").Value = "QMFJ1R"
$ws.Range("C144").Value = "NZXFG6"
$ws.Range("C145").Value = "QHZGVO"
$ws.Range("C146").Value = "OJ1F7F"
$ws.Range("C147").Value = "I8GSIE"
$ws.Range("C148").Value = "9APL6E"
$ws.Range("C149").Value = "9NVP1Y"
$ws.Range("C150").Value = "YY2891"
$ws.Range("C151").Value = "OQR6U2"
$ws.Range("C152").Value = "LBEUT6"
$ws.Range("C153").Value = "LRQD8Y"
$ws.Range("C154").Value = "DGCI0N"
$ws.Range("C155").Value = "MZE2OA"
$ws.Range("C156").Value = "XNPFD2"
$ws.Range("C157").Value = "YQJ84Q"
$ws.Range("C158").Value = "9B14QT"
$ws.Range("C159").Value = "RIXP1F"
$ws.Range("C160").Value = "MK17NH"
$ws.Range("C161").Value = "3DGGCM"
$ws.Range("C162").Value = "4W3RPI"
$ws.Range("C163").Value = "RQXT7F"
$ws.Range("C164").Value = "FEV23X"
$ws.Range("C165").Value = "2WXYMR"
$ws.Range("C166").Value = "JZ95J7"
$ws.Range("C167").Value = "QM25VX"
$ws.Range("C168").Value = "UFBJJF"
$ws.Range("C169").Value = "5GDO6Z"
$ws.Range("C170").Value = "V9YF8Q"
$ws.Range("C171").Value = "4ZS89I"
$ws.Range("C172").Value = "MEUVQQ"
$ws.Range("C173").Value = "SQZY5Z"
$ws.Range("C174").Value = "QLBKIB"
$ws.Range("C175").Value = "MXBR8J"
$ws.Range("C176").Value = "PCEEHE"
$ws.Range("C177").Value = "2B1XTT"
$ws.Range("C178").Value = "6G2FLT"
$ws.Range("C179").Value = "G7JD8A"
$ws.Range("C180").Value = "JE4NJ1"
$ws.Range("C181").Value = "2RG73Z"
$ws.Range("C182").Value = "6LOM2A"
$ws.Range("C183").Value = "K91K19"
$ws.Range("C184").Value = "II0DOH"
$ws.Range("C185").Value = "BHXKI4"
$ws.Range("C186").Value = "H0G6KN"
$ws.Range("C187").Value = "8AMSRN"
$ws.Range("C188").Value = "FX73RD"
$ws.Range("C189").Value = "VXQKY9"
$ws.Range("C190").Value = "ZY6TN9"
$ws.Range("C191").Value = "04A5WJ"
$ws.Range("C192").Value = "5OMMLT"
$ws.Range("C193").Value = "C7DEO6"
$ws.Range("C194").Value = "YL8ZAN"
$ws.Range("C195").Value = "5RJQN7"
$ws.Range("C196").Value = "4JVIRD"
$ws.Range("C197").Value = "1GIB5A"
$ws.Range("C198").Value = "340WWO"
$ws.Range("C199").Value = "HQTZY1"
$ws.Range("C200").Value = "1KAKBA"
$ws.Range("C201").Value = "DOZ3JA"
$ws.Range("C202").Value = "32Y0AS"
$ws.Range("C203").Value = "N8CQP8"
$ws.Range("C204").Value = "0MDJDC"
$ws.Range("C205").Value = "ZCYY51"
$ws.Range("C206").Value = "CPMNN2"
$ws.Range("C207").Value = "O0Z7EI"
$ws.Range("C208").Value = "2EAMS5"
$ws.Range("C209").Value = "BTGWPX"
$ws.Range("C210").Value = "AM08IR"
$ws.Range("C211").Value = "EQA060"
$ws.Range("C212").Value = "YD6N43"
$ws.Range("C213").Value = "PJQTUG"
$ws.Range("C214").Value = "4ABLKB"
$ws.Range("C215").Value = "VQOIG0"
$ws.Range("C216").Value = "G2IL9A"
$ws.Range("C217").Value = "9AFGNX"
$ws.Range("C218").Value = "S9QK8V"
$ws.Range("C219").Value = "A3H0PU"
$ws.Range("C220").Value = "RCBUOI"
$ws.Range("C221").Value = "FHMJZF"
$ws.Range("C222").Value = "3M0XMB"
$ws.Range("C223").Value = "I16498"
$ws.Range("C224").Value = "2KB3DA"
$ws.Range("C225").Value = "PB0J51"
$ws.Range("C226").Value = "M4AC5I"
$ws.Range("C227").Value = "EDUWUZ"
$ws.Range("C228").Value = "5QT24T"
$ws.Range("C229").Value = "8QWHPP"
$ws.Range("C230").Value = "69OUUM"
$ws.Range("C231").Value = "CXIRKO"
$ws.Range("C232").Value = "MO5CC9"
$ws.Range("C233").Value = "UD20W4"
$ws.Range("C234").Value = "DUUFS5"
$ws.Range("C235").Value = "0OYXS5"
$ws.Range("C236").Value = "XBV5PE"
$ws.Range("C237").Value = "82YX99"
$ws.Range("C238").Value = "IOSZMN"
$ws.Range("C239").Value = "ELWJ6S"
$ws.Range("C240").Value = "H1L8BQ"
$ws.Range("C241").Value = "R1EXG8"
$ws.Range("C242").Value = "926VCD"
$ws.Range("C243").Value = "4PUE8X"
$ws.Range("C244").Value = "TEF801"
$ws.Range("C245").Value = "Z2IH1D"
$ws.Range("C246").Value = "EP9B1V"
$ws.Range("C247").Value = "M4RA1V"
$ws.Range("C248").Value = "U4KEQU"
$ws.Range("C249").Value = "6NBL3K"
$ws.Range("C250").Value = "K3H96R"
$ws.Range("C251").Value = "MGOOPQ"
$ws.Range("C252").Value = "QTSS4K"
$ws.Range("C253").Value = "I7N1QK"
$ws.Range("C254").Value = "P8GO54"
$ws.Range("C255").Value = "4L0S57"
$ws.Range("C256").Value = "W6PUHB"
$ws.Range("C257").Value = "Z2ERPF"
$ws.Range("C258").Value = "VP85AP"
$ws.Range("C259").Value = "OLU7G9"
$ws.Range("C260").Value = "TI5NFZ"
$ws.Range("C261").Value = "CPMBET"
$ws.Range("C262").Value = "N30S1G"
$ws.Range("C263").Value = "X3F1VH"
$ws.Range("C264").Value = "FDEY5N"
$ws.Range("C265").Value = "AKUYSM"
$ws.Range("C266").Value = "H944DQ"
$ws.Range("C267").Value = "XCNQOB"
$ws.Range("C268").Value = "2TQT87"
$ws.Range("C269").Value = "WUNQVB"
$ws.Range("C270").Value = "VQEAAO"
$ws.Range("C271").Value = "YA1YFH"
$ws.Range("C272").Value = "UVZ588"
$ws.Range("C273").Value = "CGZ5C4"
$ws.Range("C274").Value = "52V48D"
$ws.Range("C275").Value = "IXMW0E"
$ws.Range("C276").Value = "Q9TANW"
$ws.Range("C277").Value = "GJHZ2Z"
$ws.Range("C278").Value = "X2AL4J"
$ws.Range("C279").Value = "WM9VJJ"
$ws.Range("C280").Value = "APMZ2P"
$ws.Range("C281").Value = "9WYGPR"
$ws.Range("C282").Value = "4JI1II"
$ws.Range("C283").Value = "GVEXYI"
$ws.Range("C284").Value = "F4HVT5"
$ws.Range("C285").Value = "PBWGDD"
$ws.Range("C286").Value = "1OZ6HW"
$ws.Range("C287").Value = "SHV1KU"
$ws.Range("C288").Value = "MBYW8F"
$ws.Range("C289").Value = "K969QC"
$ws.Range("C290").Value = "9D5BVL"
$ws.Range("C291").Value = "1FKKHA"
$ws.Range("C292").Value = "RKG4FP"
$ws.Range("C293").Value = "9RTSF0"
$ws.Range("C294").Value = "H430SF"
$ws.Range("C295").Value = "7VNREX"
$ws.Range("C296").Value = "M971RN"
$ws.Range("C297").Value = "GJBEYJ"
$ws.Range("C298").Value = "ZFRZGH"
$ws.Range("C299").Value = "VH6J4U"
$ws.Range("C300").Value = "3QUND7"
$ws.Range("C301").Value = "Z15IBB"
